$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (not numeric-looking, safe to set directly)
$ws.Range("D2").Value = '66.580.10'
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").Value = '3.187.34'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("E5").Value = '  -4.09%  '
$ws.Range("E6").Value = '  -15.10%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.176.25'
$ws.Range("E8").Value = '  -4.30%  '
$ws.Range("E9").Value = '  -9.54%  '
$ws.Range("E10").Value = '  -13.43%  '
$ws.Range("E11").Value = '  -3.82%  '
$ws.Range("E12").Value = '  -10.24%  '
$ws.Range("E13").Value = '  -10.53%  '
$ws.Range("E14").Value = '  -14.94%  '
$ws.Range("D15").Value = '3.702.83'
$ws.Range("E15").Value = '  -4.12%  '
$ws.Range("D16").Value = '66.593.35'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").Value = '3.194.05'
$ws.Range("E17").Value = '  -3.98%  '
$ws.Range("E18").Value = '  -4.67%  '
$ws.Range("E19").Value = '  -12.77%  '
$ws.Range("E20").Value = '  -10.75%  '
$ws.Range("E21").Value = '  -12.26%  '
$ws.Range("E22").Value = '  -11.42%  '
$ws.Range("E23").Value = '  -13.73%  '
$ws.Range("E24").Value = '  -9.92%  '
$ws.Range("E25").Value = '  -10.75%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -11.90%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E28").Value = '  -11.07%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E29").Value = '  -11.81%  '
$ws.Range("E30").Value = '  -9.43%  '
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("E32").Value = '  -6.79%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E35").Value = '  -17.96%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E36").Value = '  -15.91%  '
$ws.Range("E37").Value = '  -15.74%  '
$ws.Range("E38").Value = '  -7.93%  '
$ws.Range("E39").Value = '  -10.98%  '
$ws.Range("E40").Value = '  -15.24%  '
$ws.Range("E41").Value = '  -14.62%  '
$ws.Range("D42").Value = '2.796.10'
$ws.Range("E42").Value = '  -8.44%  '
$ws.Range("E43").Value = '  -13.79%  '
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("E45").Value = '  -10.88%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E46").Value = '  -6.02%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E47").Value = '  -14.13%  '
$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E48").Value = '  -11.60%  '
$ws.Range("D49").Value = '0.0₃0520'
$ws.Range("E49").Value = '  -16.07%  '
$ws.Range("E50").Value = '  -10.35%  '
$ws.Range("E51").Value = '  -21.11%  '

# Numeric-looking price strings in column D must be forced to text so
# Excel doesn't normalize/trim them (e.g. '7.20' -> 7.2, '1.00' -> 1).
# Apply a text number format, set the value, then restore the default
# 'Normal' style so no stray style index is left on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.111'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '493.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '27.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '54.03'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '488.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0410'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0800'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.117'
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.247'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '120.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.107'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.08'
$ws.Range("D51").Style = "Normal"
